$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.965.76'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '1.899.57'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7903'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9994'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3161'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.85'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07330'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08133'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7794'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.525'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').Value = '1.841.43'
$ws.Range('E15').Value = '  -3.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.259'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.81%  '
$ws.Range('D17').Value = '29.789.45'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +2.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.167'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9986'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9992'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.075.54'
$ws.Range('E24').Value = '  -3.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1615'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.499'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.047'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.447'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.549'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.504'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05633'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.106'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.251'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7573'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9995'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.670'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01940'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.801'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('D41').Value = '1.148.10'
$ws.Range('E41').Value = '  +12.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4481'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '74.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.980'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8590'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.56%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.908'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.38%  '
$ws.Range('B47').Value = 'SynthetixNetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.178'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9989'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.29'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.834'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.557'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.87%  '
